$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3000
$ws.Range("J43").Value = 3000
$ws.Range("L43").Value = 3000
$ws.Range("N43").Value = -3138

$ws.Range("H62").Value = 2351
$ws.Range("I62").Value = 2351
$ws.Range("K62").Value = 2351
$ws.Range("M62").Value = -1727

$ws.Range("H65").Value = 2351
$ws.Range("I65").Value = 2351
$ws.Range("K65").Value = 11755
$ws.Range("M65").Value = -8635

$ws.Range("H101").Value = 999.5
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 999.5
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 2998.5
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = -6242.5

$ws.Range("H137").Value = 3840.4666
$ws.Range("I137").Value = 2003.909
$ws.Range("K137").Value = 6011.727000000001
$ws.Range("M137").Value = -3461.727000000001

$ws.Range("H138").Value = 2479.8667
$ws.Range("I138").Value = 1597
$ws.Range("K138").Value = 4791
$ws.Range("M138").Value = 349

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2066941
$ws.Range("I32").Value = 1003311.8
$ws.Range("K32").Value = 1003311.8
$ws.Range("M32").Value = -1003024.8

$ws.Range("H61").Value = 2603.121
$ws.Range("I61").Value = 2283.913
$ws.Range("K61").Value = 2283.913
$ws.Range("M61").Value = -2071.913

$ws.Range("H74").Value = 267177.94
$ws.Range("I74").Value = 429503.3
$ws.Range("K74").Value = 429503.3
$ws.Range("M74").Value = -428629.3

$ws.Range("H77").Value = 267177.94
$ws.Range("I77").Value = 429503.3
$ws.Range("K77").Value = 2147516.5
$ws.Range("M77").Value = -2143148.5

$ws.Range("H97").Value = 846.8214
$ws.Range("I97").Value = 861.64
$ws.Range("K97").Value = 861.64
$ws.Range("M97").Value = -365.64

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H132").Value = 5096.0625
$ws.Range("I132").Value = 2966.9285
$ws.Range("K132").Value = 8900.7855
$ws.Range("M132").Value = -6370.7855

$ws.Range("H136").Value = 2603.121
$ws.Range("I136").Value = 2283.913
$ws.Range("K136").Value = 6851.739
$ws.Range("M136").Value = -4301.739

$ws.Range("H139").Value = 69998.91
$ws.Range("J139").Value = 69998.91
$ws.Range("L139").Value = 69998.91
$ws.Range("N139").Value = -80278.91

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2750
$ws.Range("I5").Value = 500
$ws.Range("J5").Value = 5000
$ws.Range("K5").Value = 500
$ws.Range("L5").Value = 5000
$ws.Range("M5").Value = -387
$ws.Range("N5").Value = -5226

$ws.Range("H94").Value = 68969620
$ws.Range("J94").Value = 1210.5714
$ws.Range("L94").Value = 1210.5714
$ws.Range("N94").Value = -2112.5714

$ws.Range("H105").Value = 12382737
$ws.Range("I105").Value = 589593.1
$ws.Range("K105").Value = 589593.1
$ws.Range("M105").Value = -587846.1

$ws.Range("H134").Value = 2353
$ws.Range("I134").Value = 1588.4166
$ws.Range("K134").Value = 4765.2498
$ws.Range("M134").Value = -2230.2498

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 25000
$ws.Range("J17").Value = 25000
$ws.Range("L17").Value = 25000
$ws.Range("N17").Value = -25348

$ws.Range("H31").Value = 4036602
$ws.Range("I31").Value = 2622.5
$ws.Range("K31").Value = 2622.5
$ws.Range("M31").Value = -2327.5

$ws.Range("H34").Value = 4036602
$ws.Range("I34").Value = 2622.5
$ws.Range("K34").Value = 2622.5
$ws.Range("M34").Value = -2420.5

$ws.Range("H105").Value = 1699.3
$ws.Range("I105").Value = 1699.3
$ws.Range("K105").Value = 1699.3
$ws.Range("M105").Value = 47.70000000000005

$ws.Range("H134").Value = 6903.0625
$ws.Range("I134").Value = 7032.0713
$ws.Range("K134").Value = 21096.2139
$ws.Range("M134").Value = -18561.2139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 517.8
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()

$ws.Range("H115").Value = 204831.8
$ws.Range("J115").Value = 292142.84
$ws.Range("L115").Value = 876428.52
$ws.Range("N115").Value = -878778.52

$ws.Range("H122").Value = 1316.6
$ws.Range("I122").Value = 506.4
$ws.Range("J122").Value = 1586.6666
$ws.Range("K122").Value = 4557.599999999999
$ws.Range("L122").Value = 14279.9994
$ws.Range("M122").Value = -2107.599999999999
$ws.Range("N122").Value = -19179.9994

$ws.Range("H139").Value = 9235.040000000001
$ws.Range("I139").Value = 13430.667
$ws.Range("K139").Value = 40292.001
$ws.Range("M139").Value = -35152.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5491.6665
$ws.Range("J70").Value = 6500.125
$ws.Range("L70").Value = 6500.125
$ws.Range("N70").Value = -7040.125

$ws.Range("H73").Value = 5491.6665
$ws.Range("J73").Value = 6500.125
$ws.Range("L73").Value = 6500.125
$ws.Range("N73").Value = -8372.125

$ws.Range("H102").Value = 5109.5864
$ws.Range("I102").Value = 1439.9166
$ws.Range("K102").Value = 1439.9166
$ws.Range("M102").Value = 182.0834

$ws.Range("H132").Value = 2832.3333
$ws.Range("I132").Value = 2622.7646
$ws.Range("K132").Value = 7868.293799999999
$ws.Range("M132").Value = -5338.293799999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 35000
$ws.Range("I11").Value = 35000
$ws.Range("K11").Value = 35000
$ws.Range("M11").Value = -34860

$ws.Range("H21").Value = 4750
$ws.Range("J21").Value = 4750
$ws.Range("L21").Value = 4750
$ws.Range("N21").Value = -5098

$ws.Range("H22").Value = 3225.25
$ws.Range("I22").Value = 3225.25
$ws.Range("K22").Value = 3225.25
$ws.Range("M22").Value = -2930.25

$ws.Range("H27").Value = 3225.25
$ws.Range("I27").Value = 3225.25
$ws.Range("K27").Value = 3225.25
$ws.Range("M27").Value = -3118.25

$ws.Range("H132").Value = 6020.909
$ws.Range("I132").Value = 6149.5
$ws.Range("K132").Value = 18448.5
$ws.Range("M132").Value = -15918.5

$ws.Range("H136").Value = 3219.087
$ws.Range("I136").Value = 2501.95
$ws.Range("K136").Value = 7505.849999999999
$ws.Range("M136").Value = -4955.849999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1137.7778
$ws.Range("J96").Value = 734.5
$ws.Range("L96").Value = 734.5
$ws.Range("N96").Value = -3480.5

$ws.Range("H117").Value = 100409
$ws.Range("J117").Value = 100409
$ws.Range("L117").Value = 100409
$ws.Range("N117").Value = -109587

$ws.Range("H126").Value = 2415.6667
$ws.Range("I126").Value = 2415.6667
$ws.Range("K126").Value = 7247.000100000001
$ws.Range("M126").Value = -4777.000100000001

$ws.Range("H132").Value = 2210.6428
$ws.Range("I132").Value = 1651.421
$ws.Range("K132").Value = 4954.263
$ws.Range("M132").Value = -2424.263

$ws.Range("H136").Value = 2799.0908
$ws.Range("I136").Value = 2382.72
$ws.Range("J136").Value = 4100.25
$ws.Range("K136").Value = 7148.16
$ws.Range("L136").Value = 12300.75
$ws.Range("M136").Value = -4598.16
$ws.Range("N136").Value = -17400.75
